$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 3 values (new movie record) ---
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = "ahihitesst"
$ws.Range("D3").Value = "Âu Mỹ"
$ws.Range("F3").Value = "https://www.youtube.com/watch?v=gPMaf2aw0xM"
$ws.Range("G3").Value = "/Images/fast.png"
$ws.Range("H3").Value = "Một bộ phim hài phiêu lưu theo chân một ca sĩ nổi tiếng một thời trở về quê hương và tham gia một cuộc thi khoan giếng nước với hy vọng giành được một giải thưởng tiền mặt lớn. Nhưng những gì họ không nhận ra là một con thú bí ẩn ẩn nấp dưới lòng đất."
$ws.Range("I3").Value = "144p"

# --- Clear old formatting on the row so new styling starts fresh ---
$ws.Range("A3:L3").ClearFormats()

# --- Date/time cells ---
$ws.Range("J3").Value = 44904.847222222219
$ws.Range("K3").Value = 44904.930555555555

# --- Alignment: vertical center + wrap text across the whole row ---
$ws.Range("A3:L3").VerticalAlignment = -4108
$ws.Range("A3:L3").WrapText = $true

# --- Number format for the date/time cells ---
$ws.Range("J3:K3").NumberFormat = "mm:ss.0"

# --- Row height ---
$ws.Rows(3).RowHeight = 138

# --- Column widths ---
$ws.Columns("A").ColumnWidth = 19.140625
$ws.Columns("C").ColumnWidth = 16

# --- Selection ---
$ws.Range("D8").Select()
